# The header row (row 1) was missing the "XBP1" gene label between
# "WAR1" (column X) and "YAP1" (column Y). Insert it at column X (24th
# column), shifting the existing headers from column X through AX one
# column to the right (row 1 only - the data grid below is untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$insertCol = 24                               # column X
$lastCol = $ws.UsedRange.Columns.Count        # column AX (50) before the edit

# Shift existing header cells right by one column, working from the
# right-hand side so values are not clobbered before they are read.
for ($c = $lastCol; $c -ge $insertCol; $c--) {
    $src = $ws.Cells.Item(1, $c)
    $dst = $ws.Cells.Item(1, $c + 1)
    $dst.Value = $src.Value()
}

$ws.Cells.Item(1, $insertCol).Value = "XBP1"

# Reset the view: scroll back to the top-left corner and select the
# (now one-column-wider) header range.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B1:AY1").Select()
$ws.PageSetup.Orientation = 1   # xlPortrait
